$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. NumberFormat is forced to Text ("@") before assignment
# so numeric-looking strings (e.g. "20.90", "2.880") keep their exact
# textual representation instead of being coerced to a Double and losing
# trailing zeros, matching the source workbook where these columns are
# stored as inline strings.
$cells  = @(
    "D2","E2","D3","E3","E4","D5","E5","D6","E6","E7",
    "D8","E8","D9","E9","D10","E10","D11","E11","B12","C12",
    "D12","E12","B13","C13","D13","E13","D14","E14","D15","E15",
    "D16","E16","D17","E17","D18","E18","D19","E19","D21","E21",
    "D22","E22","D23","E23","E24","D25","E25","D26","E26","D27",
    "E27","D28","E28","D29","E29","D30","E30","D31","E31","D32",
    "E32","D33","E33","D34","E34","D35","E35","D36","E36","E37",
    "D38","E38","D39","E39","D40","E40","D41","E41","D42","E42",
    "D43","E43","D44","E44","D45","E45","D46","E46","D47","E47",
    "D48","E48","E49","D50","E50","E51"
)
$values = @(
    '27.398.92',
    '  -1.35%  ',
    '1.711.04',
    '  -1.51%  ',
    '  +0.07%  ',
    '224.32',
    '  -1.30%  ',
    '0.5327',
    '  -2.37%  ',
    '  +0.13%  ',
    '0.2665',
    '  -3.32%  ',
    '0.06605',
    '  -1.78%  ',
    '20.90',
    '  -4.83%  ',
    '0.07629',
    '  -1.90%  ',
    'Polkadot',
    'https://coinranking.com/coin/25W7FG7om+polkadot-dot',
    '4.562',
    '  -2.59%  ',
    'WrappedEther',
    'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',
    '1.716.96',
    '  -1.13%  ',
    '1.950.18',
    '  -1.32%  ',
    '0.5771',
    '  -3.42%  ',
    '0.0₅8190',
    '  -2.77%  ',
    '67.92',
    '  -1.80%  ',
    '27.401.04',
    '  -1.38%  ',
    '216.76',
    '  -4.13%  ',
    '4.671',
    '  -3.23%  ',
    '10.47',
    '  -3.84%  ',
    '5.983',
    '  -3.85%  ',
    '  +0.05%  ',
    '142.31',
    '  -3.04%  ',
    '1.732',
    '  +1.85%  ',
    '0.1217',
    '  -2.38%  ',
    '7.283',
    '  -2.26%  ',
    '16.29',
    '  -5.10%  ',
    '0.05414',
    '  -4.27%  ',
    '1.291',
    '  -1.53%  ',
    '3.507',
    '  -4.87%  ',
    '3.431',
    '  -2.25%  ',
    '1.648',
    '  -1.57%  ',
    '2.880',
    '  +0.73%  ',
    '0.9504',
    '  -2.63%  ',
    '  -1.32%  ',
    '0.5868',
    '  -1.26%  ',
    '0.01634',
    '  -2.00%  ',
    '5.864',
    '  -0.17%  ',
    '1.046.08',
    '  -0.15%  ',
    '1.004',
    '  +0.10%  ',
    '0.8424',
    '  -0.60%  ',
    '100.95',
    '  -0.98%  ',
    '1.855.29',
    '  -1.40%  ',
    '0.0₈119',
    '  +2.22%  ',
    '58.02',
    '  -1.98%  ',
    '0.4513',
    '  +1.73%  ',
    '  +0.50%  ',
    '8.127',
    '  -1.56%  ',
    '  -1.32%  '
)

for ($i = 0; $i -lt $cells.Length; $i++) {
    $rng = $ws.Range($cells[$i])
    $rng.NumberFormat = "@"
    $rng.Value = $values[$i]
}

